# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Rebuilds the worker/period data block (rows 16-27) of the "Estado de
# Cuenta" sheet with the refreshed figures pulled from the updated source
# database. Column B (doc type), C (doc #), D (worker name), E (period),
# F (valor mora) and G (salario basico) are rewritten row by row; the
# existing cell formatting/styles are left untouched since only the
# underlying values changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $tipoDoc, $numDoc, $nombre, $periodo, $valorMora, $salario) {
    $ws.Cells.Item($r, 2).Value = $tipoDoc
    $ws.Cells.Item($r, 3).Value = $numDoc
    $ws.Cells.Item($r, 4).Value = $nombre
    $ws.Cells.Item($r, 5).Value = $periodo
    $ws.Cells.Item($r, 6).Value = $valorMora
    $ws.Cells.Item($r, 7).Value = $salario
}

Set-Row 16 "CC" "1002203955" "ADRIANA KARINA ALVEAR MARRUGO" "2110" 36341 908526
Set-Row 17 "CC" "1002203955" "ADRIANA KARINA ALVEAR MARRUGO" "1804" 1093  781242
Set-Row 18 "CC" "1002203955" "ADRIANA KARINA ALVEAR MARRUGO" "1803" 32800 781242
Set-Row 19 "CC" "1002203955" "ADRIANA KARINA ALVEAR MARRUGO" "1802" 32800 781242
Set-Row 20 "CC" "1002203955" "ADRIANA KARINA ALVEAR MARRUGO" "1801" 32800 781242
Set-Row 21 "CC" "1002203955" "ADRIANA KARINA ALVEAR MARRUGO" "1712" 32800 781242
Set-Row 22 "CC" "1002203955" "ADRIANA KARINA ALVEAR MARRUGO" "1711" 32800 781242
Set-Row 23 "CC" "1002203955" "ADRIANA KARINA ALVEAR MARRUGO" "1710" 32800 781242
Set-Row 24 "CC" "1051442443" "JOHANA MEDINA MARRUGO" "1809" 31249 1300000
Set-Row 25 "CC" "1051442443" "JOHANA MEDINA MARRUGO" "1808" 31249 1300000
Set-Row 26 "CC" "1051442443" "JOHANA MEDINA MARRUGO" "1807" 31249 1300000
Set-Row 27 "CC" "1051442443" "JOHANA MEDINA MARRUGO" "1806" 31249 1300000
